# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets
# to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")
$sheetExpo.Range("F2").Value = 328
$sheetExpo.Range("F3").Value = 93
$sheetExpo.Range("F4").Value = 489
$sheetExpo.Range("F5").Value = 4852
$sheetExpo.Range("F9").Value = 742
$sheetExpo.Range("F10").Value = 222

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 328
$sheetAll.Range("F3").Value = 93
$sheetAll.Range("F4").Value = 489
$sheetAll.Range("F5").Value = 4852
$sheetAll.Range("F9").Value = 742
$sheetAll.Range("F11").Value = 222
